$d = $word.ActiveDocument

# --- Paragraph: "All machine code is separated into folders ..." ---
# labled -> labeled
$null = $d.Content.Find.Execute("labled", $true, $false, $false, $false, $false, $true, 1, $false, "labeled", 2)
# non connection -> non-connection
$null = $d.Content.Find.Execute("non connection", $true, $false, $false, $false, $false, $true, 1, $false, "non-connection", 2)
# datablocks (plural) -> data blocks   (must run before singular replacement)
$null = $d.Content.Find.Execute("datablocks", $true, $false, $false, $false, $false, $true, 1, $false, "data blocks", 2)
# datablock (singular) -> data block
$null = $d.Content.Find.Execute("datablock", $true, $false, $false, $false, $false, $true, 1, $false, "data block", 2)

# --- Paragraph: "Since all machines talk back through one connection ..." ---
# "in order to send" -> "to send"
$null = $d.Content.Find.Execute("in order to send", $true, $false, $false, $false, $false, $true, 1, $false, "to send", 2)
# "that tie" -> "that, tie"
$null = $d.Content.Find.Execute("that tie", $true, $false, $false, $false, $false, $true, 1, $false, "that, tie", 2)

# --- Paragraph: "All IO is setup so that from 10.0 until 19.7 ..." ---
# The "_GoBack" bookmark previously sat in the middle of the sentence (between
# "all in" and "puts and outputs..."). Re-write the full sentence text as a
# single contiguous run (this also removes the old, mid-sentence bookmark),
# then re-create the "_GoBack" bookmark inside a brand new, separate empty
# paragraph placed right after this one.
$oldSentence = "All IO is setup so that from 10.0 until 19.7 all inputs and outputs are assigned to Machine 1 with each consecutive machine taking the next range of 10 inputs."
$null = $d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $oldSentence, 2)

# Find the paragraph that now contains this sentence and split a new, empty
# paragraph off after it.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $oldSentence) {
        $para.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i + 1)

        # Insert a temporary placeholder character so the new paragraph has a
        # real, unambiguous text range to attach the bookmark to.
        $newPara.Range.InsertAfter("X")
        $placeholderRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
        $d.Bookmarks.Add("_GoBack", $placeholderRange)

        # Remove the placeholder text again, leaving just the bookmark tags.
        $placeholderRange2 = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
        $placeholderRange2.Text = ""

        break
    }
}
